# Auto-generated edit script: apply cell-value corrections described in the commit diff.
# The workbook's underlying raw log data changed (late-arrival / threshold records were
# removed for several 進貨控場 (incoming-goods-control) team members), which cascades into
# recalculated counts/ratios across team_df, team_df_day, productivity_tl and
# productivity_team_function. Apply the resulting target values directly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("team_df")
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 29
$ws.Range("U6").Value = 0.03448275862068965
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 27
$ws.Range("U7").Value = 0.03703703703703703
$ws.Range("S15").Value = 4
$ws.Range("T15").Value = 31
$ws.Range("U15").Value = 0.1290322580645161
$ws.Range("T30").Value = 28
$ws.Range("U30").Value = 0.1785714285714286
$ws.Range("T31").Value = 28
$ws.Range("U31").Value = 0.1785714285714286
$ws.Range("T118").Value = 30
$ws.Range("U118").Value = 0.06666666666666667
$ws.Range("T119").Value = 30
$ws.Range("U119").Value = 0.06666666666666667
$ws.Range("T121").Value = 32
$ws.Range("U121").Value = 0.0625
$ws.Range("T123").Value = 37
$ws.Range("U123").Value = 0.1351351351351351
$ws.Range("T127").Value = 35
$ws.Range("U127").Value = 0.1428571428571428
$ws.Range("T128").Value = 30
$ws.Range("U128").Value = 0.1666666666666667
$ws.Range("T146").Value = 41
$ws.Range("U146").Value = 0.2682926829268293
$ws.Range("T147").Value = 49
$ws.Range("U147").Value = 0.2448979591836735
$ws.Range("T156").Value = 43
$ws.Range("U156").Value = 0.09302325581395349
$ws.Range("T160").Value = 42
$ws.Range("U160").Value = 0.07142857142857142
$ws.Range("T173").Value = 38
$ws.Range("T174").Value = 43
$ws.Range("U174").Value = 0.04651162790697674
$ws.Range("T180").Value = 25
$ws.Range("U180").Value = 0.04
$ws.Range("T183").Value = 40
$ws.Range("U183").Value = 0.1
$ws.Range("T186").Value = 32
$ws.Range("U186").Value = 0.125
$ws.Range("S199").Value = 0
$ws.Range("T199").Value = 20
$ws.Range("U199").Value = 0
$ws.Range("S201").Value = 0
$ws.Range("T201").Value = 18
$ws.Range("U201").Value = 0
$ws.Range("S204").Value = 2
$ws.Range("T204").Value = 25
$ws.Range("U204").Value = 0.08
$ws.Range("S205").Value = 3
$ws.Range("T205").Value = 22
$ws.Range("U205").Value = 0.1363636363636364
$ws.Range("T240").Value = 18
$ws.Range("U240").Value = 0.1666666666666667
$ws.Range("T244").Value = 31
$ws.Range("U244").Value = 0.2258064516129032
$ws.Range("S262").Value = 2
$ws.Range("T262").Value = 25
$ws.Range("U262").Value = 0.08
$ws.Range("S281").Value = 2
$ws.Range("T281").Value = 33
$ws.Range("U281").Value = 0.06060606060606061
$ws.Range("S286").Value = 3
$ws.Range("T286").Value = 32
$ws.Range("U286").Value = 0.09375
$ws.Range("S298").Value = 0
$ws.Range("T298").Value = 32
$ws.Range("U298").Value = 0
$ws.Range("S299").Value = 0
$ws.Range("T299").Value = 27
$ws.Range("U299").Value = 0
$ws.Range("T304").Value = 36
$ws.Range("U304").Value = 0.05555555555555555
$ws.Range("T305").Value = 34
$ws.Range("U305").Value = 0.05882352941176471
$ws.Range("T320").Value = 22
$ws.Range("U320").Value = 0.2272727272727273
$ws.Range("T321").Value = 19
$ws.Range("U321").Value = 0.2631578947368421
$ws.Range("T326").Value = 34
$ws.Range("U326").Value = 0.1764705882352941
$ws.Range("T398").Value = 19
$ws.Range("U398").Value = 0.1052631578947368
$ws.Range("T418").Value = 33
$ws.Range("U418").Value = 0.1515151515151515
$ws.Range("T433").Value = 14
$ws.Range("U433").Value = 0.07142857142857142
$ws.Range("S435").Value = 2
$ws.Range("T435").Value = 27
$ws.Range("U435").Value = 0.07407407407407407
$ws.Range("T447").Value = 19
$ws.Range("U447").Value = 0.1052631578947368
$ws.Range("S448").Value = 3
$ws.Range("T448").Value = 22
$ws.Range("U448").Value = 0.1363636363636364
$ws.Range("S451").Value = 6
$ws.Range("T451").Value = 31
$ws.Range("U451").Value = 0.1935483870967742
$ws.Range("S453").Value = 8
$ws.Range("T453").Value = 29
$ws.Range("U453").Value = 0.2758620689655172
$ws.Range("S465").Value = 0
$ws.Range("T465").Value = 14
$ws.Range("U465").Value = 0
$ws.Range("S466").Value = 0
$ws.Range("T466").Value = 11
$ws.Range("U466").Value = 0
$ws.Range("T471").Value = 24
$ws.Range("U471").Value = 0.04166666666666666
$ws.Range("T472").Value = 24
$ws.Range("U472").Value = 0.04166666666666666
$ws.Range("T486").Value = 17
$ws.Range("U486").Value = 0.1176470588235294
$ws.Range("T487").Value = 18
$ws.Range("U487").Value = 0.1111111111111111
$ws.Range("T531").Value = 8
$ws.Range("U531").Value = 0.125
$ws.Range("S544").Value = 0
$ws.Range("T544").Value = 16
$ws.Range("U544").Value = 0
$ws.Range("T547").Value = 21
$ws.Range("U547").Value = 0.2380952380952381
$ws.Range("T559").Value = 26
$ws.Range("U559").Value = 0.03846153846153846
$ws.Range("T564").Value = 34
$ws.Range("T579").Value = 25
$ws.Range("S583").Value = 1
$ws.Range("T583").Value = 27
$ws.Range("U583").Value = 0.03703703703703703
$ws.Range("T598").Value = 26
$ws.Range("T599").Value = 25
$ws.Range("T606").Value = 29
$ws.Range("U606").Value = 0.103448275862069
$ws.Range("T621").Value = 22
$ws.Range("U621").Value = 0.1818181818181818
$ws.Range("T622").Value = 20
$ws.Range("U622").Value = 0.2
$ws.Range("T624").Value = 25
$ws.Range("U624").Value = 0.08
$ws.Range("T627").Value = 24
$ws.Range("U627").Value = 0.08333333333333333

$ws = $wb.Worksheets.Item("team_df_day")
$ws.Range("G14").Value = 69
$ws.Range("H14").Value = 0.1014492753623188
$ws.Range("G15").Value = 95
$ws.Range("H15").Value = 0.05263157894736842
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 58
$ws.Range("H16").Value = 0.103448275862069
$ws.Range("G17").Value = 70
$ws.Range("H17").Value = 0.2285714285714286
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 68
$ws.Range("H18").Value = 0.1470588235294118
$ws.Range("F19").Value = 5
$ws.Range("G19").Value = 65
$ws.Range("H19").Value = 0.07692307692307693
$ws.Range("G22").Value = 59
$ws.Range("H22").Value = 0.1525423728813559
$ws.Range("G23").Value = 33
$ws.Range("H23").Value = 0.1515151515151515
$ws.Range("F24").Value = 3
$ws.Range("G24").Value = 41
$ws.Range("H24").Value = 0.07317073170731707
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = 0.1304347826086956
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 61
$ws.Range("H27").Value = 0.1147540983606557
$ws.Range("F59").Value = 5
$ws.Range("G59").Value = 60
$ws.Range("H59").Value = 0.08333333333333333
$ws.Range("G60").Value = 50
$ws.Range("H60").Value = 0.16
$ws.Range("G62").Value = 65
$ws.Range("H62").Value = 0.1076923076923077
$ws.Range("G63").Value = 92
$ws.Range("H63").Value = 0.1739130434782609
$ws.Range("G64").Value = 83
$ws.Range("H64").Value = 0.07228915662650602
$ws.Range("F65").Value = 2
$ws.Range("G65").Value = 66
$ws.Range("H65").Value = 0.0303030303030303
$ws.Range("G66").Value = 55
$ws.Range("H66").Value = 0.1818181818181818
$ws.Range("F68").Value = 11
$ws.Range("G68").Value = 51
$ws.Range("H68").Value = 0.2156862745098039
$ws.Range("F69").Value = 1
$ws.Range("G69").Value = 38
$ws.Range("H69").Value = 0.02631578947368421
$ws.Range("G70").Value = 46
$ws.Range("H70").Value = 0.1739130434782609
$ws.Range("G71").Value = 60
$ws.Range("H71").Value = 0.01666666666666667
$ws.Range("F72").Value = 1
$ws.Range("G72").Value = 52
$ws.Range("H72").Value = 0.01923076923076923
$ws.Range("G73").Value = 55
$ws.Range("H73").Value = 0.05454545454545454
$ws.Range("G74").Value = 46
$ws.Range("H74").Value = 0.1304347826086956
$ws.Range("F161").Value = 1
$ws.Range("G161").Value = 27
$ws.Range("H161").Value = 0.03703703703703703
$ws.Range("G162").Value = 48
$ws.Range("H162").Value = 0.1458333333333333
$ws.Range("G163").Value = 84
$ws.Range("H163").Value = 0.1071428571428571
$ws.Range("G164").Value = 83
$ws.Range("H164").Value = 0.1686746987951807
$ws.Range("F165").Value = 3
$ws.Range("G165").Value = 42
$ws.Range("H165").Value = 0.07142857142857142
$ws.Range("F166").Value = 2
$ws.Range("G166").Value = 63
$ws.Range("H166").Value = 0.03174603174603174
$ws.Range("G167").Value = 53
$ws.Range("H167").Value = 0.2075471698113208
$ws.Range("F169").Value = 8
$ws.Range("G169").Value = 50
$ws.Range("H169").Value = 0.16
$ws.Range("F170").Value = 1
$ws.Range("G170").Value = 35
$ws.Range("H170").Value = 0.02857142857142857
$ws.Range("G171").Value = 46
$ws.Range("H171").Value = 0.1739130434782609
$ws.Range("G172").Value = 47
$ws.Range("H172").Value = 0.0425531914893617
$ws.Range("G173").Value = 45
$ws.Range("H173").Value = 0.1333333333333333

$ws = $wb.Worksheets.Item("productivity_tl")
$ws.Range("D3").Value = 0.1055092368297509
$ws.Range("D6").Value = 0.09267044997616336
$ws.Range("D15").Value = 0.1062270329133872

$ws = $wb.Worksheets.Item("productivity_team_function")
$ws.Range("I5").Value = 0.1014492753623188
$ws.Range("K5").Value = 0.05263157894736842
$ws.Range("L5").Value = 0.103448275862069
$ws.Range("N5").Value = 0.2285714285714286
$ws.Range("O5").Value = 0.1470588235294118
$ws.Range("P5").Value = 0.07692307692307693
$ws.Range("U5").Value = 0.1525423728813559
$ws.Range("V5").Value = 0.1515151515151515
$ws.Range("W5").Value = 0.07317073170731707
$ws.Range("AB5").Value = 0.1304347826086956
$ws.Range("AC5").Value = 0.1147540983606557
$ws.Range("D9").Value = 0.08333333333333333
$ws.Range("E9").Value = 0.16
$ws.Range("I9").Value = 0.1076923076923077
$ws.Range("J9").Value = 0.1739130434782609
$ws.Range("K9").Value = 0.07228915662650602
$ws.Range("Q9").Value = 0.0303030303030303
$ws.Range("R9").Value = 0.1818181818181818
$ws.Range("X9").Value = 0.2156862745098039
$ws.Range("Y9").Value = 0.02631578947368421
$ws.Range("Z9").Value = 0.1739130434782609
$ws.Range("AD9").Value = 0.01666666666666667
$ws.Range("AE9").Value = 0.01923076923076923
$ws.Range("AF9").Value = 0.05454545454545454
$ws.Range("AG9").Value = 0.1304347826086956
$ws.Range("D18").Value = 0.03703703703703703
$ws.Range("E18").Value = 0.1458333333333333
$ws.Range("I18").Value = 0.1071428571428571
$ws.Range("J18").Value = 0.1686746987951807
$ws.Range("L18").Value = 0.07142857142857142
$ws.Range("Q18").Value = 0.03174603174603174
$ws.Range("R18").Value = 0.2075471698113208
$ws.Range("X18").Value = 0.16
$ws.Range("Y18").Value = 0.02857142857142857
$ws.Range("Z18").Value = 0.1739130434782609
$ws.Range("AF18").Value = 0.0425531914893617
$ws.Range("AG18").Value = 0.1333333333333333
